$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 24,14
$data[0,0] = 1.057073508126791
$data[0,1] = 0.1760460016094498
$data[0,2] = 0.58421196770405
$data[0,3] = 0.2146193540449808
$data[0,4] = 0
$data[0,5] = 0.002480050077520375
$data[0,6] = 0
$data[0,7] = 0
$data[0,8] = 0.09565057122241427
$data[0,9] = 0
$data[0,10] = 0
$data[0,11] = 0.4984823096152411
$data[0,12] = 1.493658391883164
$data[0,13] = 4.615361159481807
$data[1,0] = 0.9737045013519605
$data[1,1] = 0.15678593737195
$data[1,2] = 0.5803722276089189
$data[1,3] = 0.2142104062548107
$data[1,4] = 0
$data[1,5] = 0.002483661301416862
$data[1,6] = 0
$data[1,7] = 0
$data[1,8] = 0.096044811201903
$data[1,9] = 0
$data[1,10] = 0
$data[1,11] = 0.4757391542780098
$data[1,12] = 1.512240688699801
$data[1,13] = 4.607621046267212
$data[2,0] = 0.9228478474239807
$data[2,1] = 0.1449711080410623
$data[2,2] = 0.5782943428948784
$data[2,3] = 0.2140577605890819
$data[2,4] = 0
$data[2,5] = 0.002485997403705619
$data[2,6] = 0
$data[2,7] = 0
$data[2,8] = 0.09632602706154714
$data[2,9] = 0
$data[2,10] = 0
$data[2,11] = 0.4619962172647831
$data[2,12] = 1.524241352484939
$data[2,13] = 4.605706632628142
$data[3,0] = 0.9022075286146958
$data[3,1] = 0.1401593608371456
$data[3,2] = 0.577517961496639
$data[3,3] = 0.2140203181345548
$data[3,4] = 0
$data[3,5] = 0.002486979350455506
$data[3,6] = 0
$data[3,7] = 0
$data[3,8] = 0.09645046953352931
$data[3,9] = 0
$data[3,10] = 0
$data[3,11] = 0.4564516985674842
$data[3,12] = 1.529280282877249
$data[3,13] = 4.605639078046067
$data[4,0] = 0.8987853225280844
$data[4,1] = 0.139360552806977
$data[4,2] = 0.5773932953426595
$data[4,3] = 0.2140155965843284
$data[4,4] = 0
$data[4,5] = 0.002487144214803882
$data[4,6] = 0
$data[4,7] = 0
$data[4,8] = 0.09647172770931078
$data[4,9] = 0
$data[4,10] = 0
$data[4,11] = 0.4555344126372773
$data[4,12] = 1.530125964445517
$data[4,13] = 4.605670868054119
$data[5,0] = 0.9225691431802829
$data[5,1] = 0.1449062032351094
$data[5,2] = 0.5782835873590955
$data[5,3] = 0.2140571553588906
$data[5,4] = 0
$data[5,5] = 0.002486010525126027
$data[5,6] = 0
$data[5,7] = 0
$data[5,8] = 0.09632766547602145
$data[5,9] = 0
$data[5,10] = 0
$data[5,11] = 0.4619212156484167
$data[5,12] = 1.524308707826384
$data[5,13] = 4.605702837763232
$data[6,0] = 1.028259266080909
$data[6,1] = 0.1694029432250375
$data[6,2] = 0.5828299838237427
$data[6,3] = 0.2144579183267332
$data[6,4] = 0
$data[6,5] = 0.002481270628074128
$data[6,6] = 0
$data[6,7] = 0
$data[6,8] = 0.09577837920745225
$data[6,9] = 0
$data[6,10] = 0
$data[6,11] = 0.4905945907683389
$data[6,12] = 1.49994285622868
$data[6,13] = 4.61210267262075
$data[7,0] = 1.238137140164156
$data[7,1] = 0.217523841066253
$data[7,2] = 0.5939647750066683
$data[7,3] = 0.2160250249359983
$data[7,4] = 0
$data[7,5] = 0.002472914014870512
$data[7,6] = 0
$data[7,7] = 0
$data[7,8] = 0.09501193539940544
$data[7,9] = 0
$data[7,10] = 0
$data[7,11] = 0.548576817761095
$data[7,12] = 1.456853523719051
$data[7,13] = 4.647227044471975
$data[8,0] = 1.39392632699014
$data[8,1] = 0.2529272715508455
$data[8,2] = 0.603499980394929
$data[8,3] = 0.2176531658428793
$data[8,4] = 0
$data[8,5] = 0.002467340429454849
$data[8,6] = 0
$data[8,7] = 0
$data[8,8] = 0.09463840736130535
$data[8,9] = 0
$data[8,10] = 0
$data[8,11] = 0.5922462195547382
$data[8,12] = 1.428058365455049
$data[8,13] = 4.6868829241144
$data[9,0] = 1.465144565588503
$data[9,1] = 0.2690437789334226
$data[9,2] = 0.6081323668762479
$data[9,3] = 0.2184975330218464
$data[9,4] = 0
$data[9,5] = 0.002464926488943894
$data[9,6] = 0
$data[9,7] = 0
$data[9,8] = 0.09450968925884595
$data[9,9] = 0
$data[9,10] = 0
$data[9,11] = 0.6123453722971846
$data[9,12] = 1.415580810039282
$data[9,13] = 4.707950642147978
$data[10,0] = 1.492162870805601
$data[10,1] = 0.2751482202019702
$data[10,2] = 0.6099289238021015
$data[10,3] = 0.2188321922668734
$data[10,4] = 0
$data[10,5] = 0.002464029768756738
$data[10,6] = 0
$data[10,7] = 0
$data[10,8] = 0.09446687477247195
$data[10,9] = 0
$data[10,10] = 0
$data[10,11] = 0.6199899394272563
$data[10,12] = 1.410945347182061
$data[10,13] = 4.716365295469018
$data[11,0] = 1.486341804346864
$data[11,1] = 0.2738334566627429
$data[11,2] = 0.6095401188192966
$data[11,3] = 0.2187594540059585
$data[11,4] = 0
$data[11,5] = 0.002464222121678745
$data[11,6] = 0
$data[11,7] = 0
$data[11,8] = 0.09447583191944631
$data[11,9] = 0
$data[11,10] = 0
$data[11,11] = 0.6183420596756548
$data[11,12] = 1.411939693411105
$data[11,13] = 4.714533602683559
$data[12,0] = 1.46736638953513
$data[12,1] = 0.2695459664272732
$data[12,2] = 0.6082793215672382
$data[12,3] = 0.2185247666974632
$data[12,4] = 0
$data[12,5] = 0.002464852367148329
$data[12,6] = 0
$data[12,7] = 0
$data[12,8] = 0.09450604805312324
$data[12,9] = 0
$data[12,10] = 0
$data[12,11] = 0.6129736258433383
$data[12,12] = 1.415197652827326
$data[12,13] = 4.708634159627138
$data[13,0] = 1.455749826984515
$data[13,1] = 0.2669199416376387
$data[13,2] = 0.6075125645863579
$data[13,3] = 0.218382956455212
$data[13,4] = 0
$data[13,5] = 0.002465240672825968
$data[13,6] = 0
$data[13,7] = 0
$data[13,8] = 0.09452532845023498
$data[13,9] = 0
$data[13,10] = 0
$data[13,11] = 0.609689657274771
$data[13,12] = 1.417204909126577
$data[13,13] = 4.705077499975971
$data[14,0] = 1.389279011604856
$data[14,1] = 0.2518742340427593
$data[14,2] = 0.6032031733011536
$data[14,3] = 0.2176000716242577
$data[14,4] = 0
$data[14,5] = 0.002467500623738913
$data[14,6] = 0
$data[14,7] = 0
$data[14,8] = 0.09464764877155929
$data[14,9] = 0
$data[14,10] = 0
$data[14,11] = 0.5909373824684465
$data[14,12] = 1.428886321593611
$data[14,13] = 4.685567145338212
$data[15,0] = 1.348590212735587
$data[15,1] = 0.2426469644389329
$data[15,2] = 0.6006349886226019
$data[15,3] = 0.2171463624851242
$data[15,4] = 0
$data[15,5] = 0.002468918090833174
$data[15,6] = 0
$data[15,7] = 0
$data[15,8] = 0.09473324351063539
$data[15,9] = 0
$data[15,10] = 0
$data[15,11] = 0.5794932204257961
$data[15,12] = 1.436211773787699
$data[15,13] = 4.674374679779589
$data[16,0] = 1.325219977190443
$data[16,1] = 0.2373407567609718
$data[16,2] = 0.5991855830699251
$data[16,3] = 0.21689516329198
$data[16,4] = 0
$data[16,5] = 0.002469744822798152
$data[16,6] = 0
$data[16,7] = 0
$data[16,8] = 0.09478635321581663
$data[16,9] = 0
$data[16,10] = 0
$data[16,11] = 0.5729328493308472
$data[16,12] = 1.440483670262493
$data[16,13] = 4.668222004446704
$data[17,0] = 1.317312885777483
$data[17,1] = 0.2355443552268071
$data[17,2] = 0.5986996052046436
$data[17,3] = 0.2168117882152778
$data[17,4] = 0
$data[17,5] = 0.00247002670811617
$data[17,6] = 0
$data[17,7] = 0
$data[17,8] = 0.09480500117463464
$data[17,9] = 0
$data[17,10] = 0
$data[17,11] = 0.5707154050882508
$data[17,12] = 1.441940102428418
$data[17,13] = 4.666187711656164
$data[18,0] = 1.352918206234676
$data[18,1] = 0.2436291130520374
$data[18,2] = 0.6009055047385061
$data[18,3] = 0.2171936502283245
$data[18,4] = 0
$data[18,5] = 0.002468766015701235
$data[18,6] = 0
$data[18,7] = 0
$data[18,8] = 0.09472373043378823
$data[18,9] = 0
$data[18,10] = 0
$data[18,11] = 0.5807091943723606
$data[18,12] = 1.43542591235008
$data[18,13] = 4.675536637504536
$data[19,0] = 1.472938590667752
$data[19,1] = 0.2708052675928343
$data[19,2] = 0.6086484984168976
$data[19,3] = 0.2185932952286613
$data[19,4] = 0
$data[19,5] = 0.002464666777806148
$data[19,6] = 0
$data[19,7] = 0
$data[19,8] = 0.09449701192227167
$data[19,9] = 0
$data[19,10] = 0
$data[19,11] = 0.6145495582790232
$data[19,12] = 1.414238280538452
$data[19,13] = 4.710355104674989
$data[20,0] = 1.551667087601857
$data[20,1] = 0.2885749505859394
$data[20,2] = 0.6139559666175955
$data[20,3] = 0.2195949852438872
$data[20,4] = 0
$data[20,5] = 0.002462088991855231
$data[20,6] = 0
$data[20,7] = 0
$data[20,8] = 0.09438339449266664
$data[20,9] = 0
$data[20,10] = 0
$data[20,11] = 0.6368611186461237
$data[20,12] = 1.400912780381592
$data[20,13] = 4.735657417171296
$data[21,0] = 1.50962205494443
$data[21,1] = 0.2790902103649842
$data[21,2] = 0.6111006775741146
$data[21,3] = 0.2190524087983015
$data[21,4] = 0
$data[21,5] = 0.002463455564034076
$data[21,6] = 0
$data[21,7] = 0
$data[21,8] = 0.09444087114964361
$data[21,9] = 0
$data[21,10] = 0
$data[21,11] = 0.6249352338162595
$data[21,12] = 1.407977052843304
$data[21,13] = 4.721919654652424
$data[22,0] = 1.350961450868567
$data[22,1] = 0.2431850877897546
$data[22,2] = 0.6007831200575708
$data[22,3] = 0.2171722413991652
$data[22,4] = 0
$data[22,5] = 0.002468834731830035
$data[22,6] = 0
$data[22,7] = 0
$data[22,8] = 0.09472801914733964
$data[22,9] = 0
$data[22,10] = 0
$data[22,11] = 0.5801593932635782
$data[22,12] = 1.435781012174395
$data[22,13] = 4.675010438004165
$data[23,0] = 1.181079704229887
$data[23,1] = 0.2044971539452547
$data[23,2] = 0.5907147729775062
$data[23,3] = 0.2155173926839105
$data[23,4] = 0
$data[23,5] = 0.00247507487020756
$data[23,6] = 0
$data[23,7] = 0
$data[23,8] = 0.09518599933001681
$data[23,9] = 0
$data[23,10] = 0
$data[23,11] = 0.5327032676816117
$data[23,12] = 1.468007810414729
$data[23,13] = 4.635299166977262

$ws.Range("B2:O25").Value = $data
